# Video Start/Stop Record Added
#
# - Rename the "movement_path" parameter to "movement_file"
# - Point the camera config at the "no_topic" variant (Cam_Basic.sdf -> Cam_Basic_no_topic.sdf)
# - Update the active selection on the "Sheet1" worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# A2 held "movement_path" -> rename to "movement_file"
$ws.Range("A2").Value = "movement_file"

# B6 held "Cam_Basic.sdf" -> rename to "Cam_Basic_no_topic.sdf"
$ws.Range("B6").Value = "Cam_Basic_no_topic.sdf"

# Move the active selection from C6 to A8
$ws.Range("A8").Select()
